# Update demo content assessment outputs:
# Replace the "files/content_assessment/" prefix with "demo_files/" in the
# Object Path column (column A) across the PLSQL, SQR, and ET worksheets.

$wb = $excel.ActiveWorkbook

$sheetsInfo = @(
    @{ Name = "PLSQL"; LastRow = 16 },
    @{ Name = "SQR";   LastRow = 6 },
    @{ Name = "ET";    LastRow = 11 }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)
    for ($r = 2; $r -le $info.LastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Text
        if ($current -ne $null) {
            $updated = $current.Replace("files/content_assessment/", "demo_files/")
            if ($updated -ne $current) {
                $cell.Formula = $updated
            }
        }
    }
}
